# Further work on report. Adds comparison between var sets.
$wb = $excel.ActiveWorkbook

# --- Add the new "Sheet3" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# --- Populate header row (row 1) ---
$ws3.Range("C1").Value = "hclust pc1-2 var"
$ws3.Range("D1").Value = "avg sil width"
$ws3.Range("E1").Value = "best k, sil width"
$ws3.Range("G1").Value = "kmeans elbow"
$ws3.Range("H1").Value = "rsq"
$ws3.Range("I1").Value = "avg sil width"
$ws3.Range("J1").Value = "best k, sil width"
$ws3.Range("L1").Value = "pam avg sil"
$ws3.Range("M1").Value = "best k"
$ws3.Range("O1").Value = "mclust k"
$ws3.Range("P1").Value = "avg sil width"

# --- Row 2 (set1) ---
$ws3.Range("A2").Value = "set1"
$ws3.Range("B2").Value = "ward.D"
$ws3.Range("C2").Value = 28.3
$ws3.Range("D2").Value = 0.093122079999999996
$ws3.Range("E2").Value = 5
$ws3.Range("G2").Value = 6
$ws3.Range("H2").Value = 0.27
$ws3.Range("I2").Value = 0.11379431
$ws3.Range("J2").Value = 6
$ws3.Range("K2").Value = "atleast one cluster has ~0 width"
$ws3.Range("L2").Value = 0.10027680999999999
$ws3.Range("M2").Value = 10
$ws3.Range("N2").Value = "lots of overlap"
$ws3.Range("O2").Value = 3
$ws3.Range("P2").Value = 0.09
$ws3.Range("Q2").Value = "1 negative width"

# --- Row 3 (set2) ---
$ws3.Range("A3").Value = "set2"
$ws3.Range("B3").Value = "ward.D2"
$ws3.Range("C3").Value = 31.87
$ws3.Range("D3").Value = 0.12016014999999999
$ws3.Range("E3").Value = 6
$ws3.Range("G3").Value = 8
$ws3.Range("H3").Value = 0.13
$ws3.Range("I3").Value = 0.13295109999999999
$ws3.Range("J3").Value = 2
$ws3.Range("L3").Value = 0.08817055
$ws3.Range("M3").Value = 9
$ws3.Range("N3").Value = "lots of overlap"
$ws3.Range("O3").Value = 3
$ws3.Range("P3").Value = 0.02
$ws3.Range("Q3").Value = "lots of overlap"

# --- Row 4 (set3) ---
$ws3.Range("A4").Value = "set3"
$ws3.Range("B4").Value = "ward.D2"
$ws3.Range("C4").Value = 25.22
$ws3.Range("D4").Value = 0.06
$ws3.Range("E4").Value = 2
$ws3.Range("G4").Value = 6
$ws3.Range("H4").Value = 0.28710000000000002
$ws3.Range("I4").Value = 0.087773470000000006
$ws3.Range("J4").Value = 7
$ws3.Range("L4").Value = 0.066646449999999996
$ws3.Range("M4").Value = 5
$ws3.Range("O4").Value = 3
$ws3.Range("P4").Value = 0.02
$ws3.Range("Q4").Value = "lots of overlap"

# --- Row 5 ---
$ws3.Range("B5").Value = "ward.D"
$ws3.Range("D5").Value = 0.33500000000000002
$ws3.Range("E5").Value = 2

# --- Selection / view state for the new sheet ---
$ws3.Range("P3").Select()

# --- Window / workbook view state ---
$win = $excel.ActiveWindow
$win.Left = 2460
$win.Top = 440
$win.Width = 23140
$win.Height = 15560
